# Auto-generated edit script: updates currentAveragePrice / Leve profit
# calculation columns (H-N) across multiple sheets, per scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2643.4285
$ws.Range("I2").Value = 155
$ws.Range("J2").Value = 3638.8
$ws.Range("K2").Value = 155
$ws.Range("L2").Value = 3638.8
$ws.Range("M2").Value = -42
$ws.Range("N2").Value = -3864.8
$ws.Range("H15").Value = 431.95456
$ws.Range("I15").Value = 431.95456
$ws.Range("K15").Value = 1295.86368
$ws.Range("M15").Value = -1126.86368
$ws.Range("H32").Value = 891.8461
$ws.Range("I32").Value = 783.1667
$ws.Range("K32").Value = 783.1667
$ws.Range("M32").Value = -457.1667
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H55").Value = 394.25
$ws.Range("I55").Value = 305
$ws.Range("J55").Value = 543
$ws.Range("K55").Value = 305
$ws.Range("L55").Value = 543
$ws.Range("M55").Value = -91
$ws.Range("N55").Value = -971
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H98").Value = 2864.6316
$ws.Range("I98").Value = 2574.0715
$ws.Range("J98").Value = 3678.2
$ws.Range("K98").Value = 2574.0715
$ws.Range("L98").Value = 3678.2
$ws.Range("M98").Value = -1076.0715
$ws.Range("N98").Value = -6674.2
$ws.Range("H111").Value = 1374.75
$ws.Range("I111").Value = 833
$ws.Range("K111").Value = 2499
$ws.Range("M111").Value = 568
$ws.Range("H112").Value = 2598.5454
$ws.Range("J112").Value = 2688.5
$ws.Range("L112").Value = 8065.5
$ws.Range("N112").Value = -10281.5
$ws.Range("H122").Value = 2864.6316
$ws.Range("I122").Value = 2574.0715
$ws.Range("J122").Value = 3678.2
$ws.Range("K122").Value = 7722.2145
$ws.Range("L122").Value = 11034.6
$ws.Range("M122").Value = -5272.2145
$ws.Range("N122").Value = -15934.6
$ws.Range("H129").Value = 3555.75
$ws.Range("I129").Value = 897.5
$ws.Range("J129").Value = 4087.4
$ws.Range("K129").Value = 2692.5
$ws.Range("L129").Value = 12262.2
$ws.Range("M129").Value = 2307.5
$ws.Range("N129").Value = -22262.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 160.8
$ws.Range("I5").Value = 201.66667
$ws.Range("J5").Value = 99.5
$ws.Range("K5").Value = 201.66667
$ws.Range("L5").Value = 99.5
$ws.Range("M5").Value = -89.66667000000001
$ws.Range("N5").Value = -323.5
$ws.Range("H63").Value = 2534.2727
$ws.Range("I63").Value = 1695.6
$ws.Range("K63").Value = 1695.6
$ws.Range("M63").Value = -1009.6
$ws.Range("H66").Value = 2534.2727
$ws.Range("I66").Value = 1695.6
$ws.Range("K66").Value = 8478
$ws.Range("M66").Value = -5046
$ws.Range("H122").Value = 2554.1428
$ws.Range("I122").Value = 2481.5833
$ws.Range("J122").Value = 2989.5
$ws.Range("K122").Value = 7444.749899999999
$ws.Range("L122").Value = 8968.5
$ws.Range("M122").Value = -4994.749899999999
$ws.Range("N122").Value = -13868.5
$ws.Range("H132").Value = 1331.5416
$ws.Range("I132").Value = 1331.5416
$ws.Range("K132").Value = 3994.6248
$ws.Range("M132").Value = -1464.6248

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 160.8
$ws.Range("I4").Value = 201.66667
$ws.Range("J4").Value = 99.5
$ws.Range("K4").Value = 201.66667
$ws.Range("L4").Value = 99.5
$ws.Range("M4").Value = -86.66667000000001
$ws.Range("N4").Value = -329.5
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H86").Value = 2650.3333
$ws.Range("I86").Value = 1671.7693
$ws.Range("K86").Value = 1671.7693
$ws.Range("M86").Value = -548.7692999999999
$ws.Range("H89").Value = 2650.3333
$ws.Range("I89").Value = 1671.7693
$ws.Range("K89").Value = 8358.8465
$ws.Range("M89").Value = -2742.8465
$ws.Range("H105").Value = 2402.25
$ws.Range("I105").Value = 2232.6667
$ws.Range("K105").Value = 2232.6667
$ws.Range("M105").Value = -485.6667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 87.76470999999999
$ws.Range("I7").Value = 102.90909
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 102.90909
$ws.Range("L7").Value = 60
$ws.Range("M7").Value = 10.09090999999999
$ws.Range("N7").Value = -286
$ws.Range("H16").Value = 2983.7058
$ws.Range("I16").Value = 1481.5454
$ws.Range("K16").Value = 1481.5454
$ws.Range("M16").Value = -1194.5454
$ws.Range("H38").Value = 11500
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 11500
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 11500
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -12254
$ws.Range("H46").Value = 11500
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 11500
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 11500
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -11922
$ws.Range("H60").Value = 21126.857
$ws.Range("J60").Value = 24082.5
$ws.Range("L60").Value = 24082.5
$ws.Range("N60").Value = -25104.5
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880
$ws.Range("H113").Value = 2983.7058
$ws.Range("I113").Value = 1481.5454
$ws.Range("K113").Value = 1481.5454
$ws.Range("M113").Value = 688.4546
$ws.Range("H122").Value = 802.25
$ws.Range("I122").Value = 567.9091
$ws.Range("K122").Value = 1703.7273
$ws.Range("M122").Value = 746.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 270.53845
$ws.Range("J12").Value = 332.7143
$ws.Range("L12").Value = 998.1428999999999
$ws.Range("N12").Value = -1344.1429
$ws.Range("H48").Value = 2000
$ws.Range("J48").Value = 2000
$ws.Range("L48").Value = 6000
$ws.Range("N48").Value = -6500
$ws.Range("H75").Value = 89.5
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 89.5
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H128").Value = 339899
$ws.Range("I128").Value = 339899
$ws.Range("K128").Value = 1019697
$ws.Range("M128").Value = -1014717
$ws.Range("H131").Value = 2109.4
$ws.Range("J131").Value = 2109.4
$ws.Range("L131").Value = 6328.200000000001
$ws.Range("N131").Value = -16408.2
$ws.Range("H139").Value = 4937.4
$ws.Range("I139").Value = 4937.4
$ws.Range("K139").Value = 14812.2
$ws.Range("M139").Value = -9672.199999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H47").Value = 10000
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11136
$ws.Range("H97").Value = 1969.25
$ws.Range("I97").Value = 2139.1428
$ws.Range("J97").Value = 780
$ws.Range("K97").Value = 2139.1428
$ws.Range("L97").Value = 780
$ws.Range("M97").Value = -1643.1428
$ws.Range("N97").Value = -1772
$ws.Range("H98").Value = 16500
$ws.Range("J98").Value = 16500
$ws.Range("L98").Value = 16500
$ws.Range("N98").Value = -22490
$ws.Range("H107").Value = 537.6
$ws.Range("I107").Value = 363
$ws.Range("J107").Value = 799.5
$ws.Range("K107").Value = 363
$ws.Range("L107").Value = 799.5
$ws.Range("M107").Value = 1557
$ws.Range("N107").Value = -4639.5
$ws.Range("H126").Value = 4000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1999
$ws.Range("I100").Value = 1999
$ws.Range("K100").Value = 1999
$ws.Range("M100").Value = -1458

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2862.5
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

